$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E10").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("E12").Value = "[-, 'MEC-2A-Des. Maq. Cad_T1']"
$ws.Range("E14").Value = "['MEC-2A-Des. Maq. Cad_T1', 'MEC-2A-Des. Maq. Cad_T1']"

$ws.Range("B18").Value = "MEC-2NB-Des. Maq. Cad"
$ws.Range("C18").Value = "['MEC-1NB-Metalografia', 'MEC-1NB-Trat. Termicos', -, -]"
$ws.Range("D18").Value = "['ELM-1NA-Des. Bas. Mec.', 'ELM-1NA-Des. Bas. Mec.']"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, 'MEC-1NA-Metalografia', -, -]"

$ws.Range("B19").Value = "MEC-2NB-Des. Maq. Cad"
$ws.Range("C19").Value = "['MEC-1NB-Metalografia', 'MEC-1NB-Trat. Termicos', -, -]"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "[-, 'MEC-1NA-Metalografia', -, -]"

$ws.Range("B20").Value = "MEC-2NB-Des. Maq. Cad"
$ws.Range("C20").Value = "['MEC-1NB-Metalografia', 'MEC-1NB-Trat. Termicos', -, -]"
$ws.Range("F20").Value = "[-, 'MEC-1NA-Metalografia', -, -]"

$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "['MEC-1NB-Metalografia', 'MEC-1NB-Trat. Termicos', -, -]"
$ws.Range("D21").Value = "['ELM-1NA-Des. Bas. Mec.', 'ELM-1NA-Des. Bas. Mec.']"
$ws.Range("F21").Value = "[-, 'MEC-1NA-Metalografia', -, -]"
